# Insert a new weekly Cilantro price record as row 46 on the "Macroferia
# Regional de Talca" sheet, pushing the existing rows 46:62 down to 47:63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 46, shifting rows 46:62 -> 47:63.
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new record.
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = 44825
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = 100112040
$ws.Range("G46").Value = "Cilantro"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 150
$ws.Range("K46").Value = 7000
$ws.Range("L46").Value = 7000
$ws.Range("M46").Value = 7000
$ws.Range("N46").Value = "`$/caja 36 atados"
$ws.Range("O46").Value = "Región del Maule"
$ws.Range("P46").Value = 194
$ws.Range("Q46").Value = 36
$ws.Range("R46").Value = "Hortaliza"
